$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Hindi MCQ question about the capital of India
$ws.Range("A7").Value = "mcq"
$ws.Range("B7").Value = "भारत की राजधानी क्या है?"
$ws.Range("C7").Value = "A:मुंबई,B:दिल्ली,C:कोलकाता,D:चेन्नई"
$ws.Range("D7").Value = "b"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = "दिल्ली भारत की राजधानी है"

# Row 8: Hindi fill-in-the-blank question about pi
$ws.Range("A8").Value = "fillblank"
$ws.Range("B8").Value = "गणित में π का मान लगभग ___ होता है।"
# The correct answer looks like a number, so force it to be stored as
# text (matching the numberStoredAsText handling already used elsewhere
# in this sheet) without leaving a stray number-format on the cell.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.14159"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "3.14159 गणित में π का लगभग मान है"
